# Update the "paquetes" (R packages) list on Hoja1.
#
# The MODELADO column (B) grows from 8 to 13 entries: two stale
# packages are replaced, and several packages that used to live in the
# SECUNDARIOS column (C) move up into MODELADO, picking up the
# underlined style already used by B2/B3. SECUNDARIOS itself shrinks
# and is renumbered to close the gap. A15 loses its underline, and the
# active selection ends on C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUnderlineStyleSingle = 2
$xlUnderlineStyleNone = -4142

# --- SECUNDARIOS (column C) -------------------------------------------
# Re-write the whole column top to bottom with the new, shorter list,
# then drop the three trailing rows that no longer have an entry.
$secundarios = @("Ckmeans.1d.dp", "tensorflow", "MASS", "DT", "rminer", "e1071", "corrgram")
for ($i = 0; $i -lt $secundarios.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $secundarios[$i]
}
$ws.Range("C9").Clear()
$ws.Range("C10").Clear()
$ws.Range("C11").Clear()

# --- MODELADO (column B) ------------------------------------------------
# B4/B5 used to hold "nortest"/"lmtest"; those packages are dropped and
# the cells now hold two packages moved up from SECUNDARIOS, underlined.
$ws.Range("B4").Value = "randomForest"
$ws.Range("B4").Font.Underline = $xlUnderlineStyleSingle

$ws.Range("B5").Value = "xgboost"
$ws.Range("B5").Font.Underline = $xlUnderlineStyleSingle

# B9 swaps its old value (tensorflow) for the new one.
$ws.Range("B9").Value = "adabag"

# New underlined entries appended below the previous end of the list.
$ws.Range("B10").Value = "xgboost"
$ws.Range("B10").Font.Underline = $xlUnderlineStyleSingle

$ws.Range("B11").Value = "rpart"
$ws.Range("B11").Font.Underline = $xlUnderlineStyleSingle

$ws.Range("B12").Value = "kernlab"

$ws.Range("B13").Value = "rpart.plot"
$ws.Range("B13").Font.Underline = $xlUnderlineStyleSingle

$ws.Range("B14").Value = "DT"
$ws.Range("B14").Font.Underline = $xlUnderlineStyleSingle

# --- BASE (column A) -----------------------------------------------------
# A15 ("reshape2") loses its underline formatting.
$ws.Range("A15").Font.Underline = $xlUnderlineStyleNone

# --- Selection -------------------------------------------------------------
$ws.Range("C9").Select()
